# Excel COM-interop edit script
# New crime data collected: refreshes the CompStat weekly report header
# (volume number + reporting week dates) and the weekly / 28-day / year-to-date
# / 2-year crime-complaint figures in the precinct comparison table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Writes a text value into a cell and re-applies the donor cells number format
# (style only, via PasteSpecial formats) so the cell keeps/gains the shared-string
# "text" formatting instead of Excels auto quote-prefix style.
function Set-TextCell {
    param($Row, $Col, $Text, $DonorRow, $DonorCol)
    $cell = $ws.Cells.Item($Row, $Col)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $ws.Cells.Item($DonorRow, $DonorCol).Copy()
    $cell.PasteSpecial(-4122)
}

# Writes a numeric value into a cell that used to hold a text placeholder, then
# re-applies the donor cells numeric style (PasteSpecial formats) so the cell
# picks up the right number format instead of keeping the old text style.
function Set-NumCell {
    param($Row, $Col, $Value, $DonorRow, $DonorCol)
    $cell = $ws.Cells.Item($Row, $Col)
    $cell.Value = $Value
    $ws.Cells.Item($DonorRow, $DonorCol).Copy()
    $cell.PasteSpecial(-4122)
}

# --- Header text: volume number and report week dates ---
$ws.Cells.Item(8, 1).Value = "Volume 29   Number  47"
$ws.Cells.Item(9, 3).Value = "Report Covering the Week  11/21/2022  Through  11/27/2022"

# --- Crime complaint figures (rows 14-30) ---
# Row 14
Set-NumCell 14 4 2 14 9
Set-NumCell 14 5 -100 14 14
Set-NumCell 14 7 2 14 9
Set-NumCell 14 8 -100 14 14
$ws.Cells.Item(14, 10).Value = 5
$ws.Cells.Item(14, 11).Value = 20
$ws.Cells.Item(14, 12).Value = 100
# Row 15
Set-TextCell 15 3 "0" 23 3
Set-TextCell 15 7 "0" 23 3
Set-TextCell 15 8 "***.*" 23 5
# Row 16
$ws.Cells.Item(16, 3).Value = 17
$ws.Cells.Item(16, 5).Value = 70
$ws.Cells.Item(16, 6).Value = 67
$ws.Cells.Item(16, 7).Value = 38
$ws.Cells.Item(16, 8).Value = 76.315789473684
$ws.Cells.Item(16, 9).Value = 594
$ws.Cells.Item(16, 10).Value = 394
$ws.Cells.Item(16, 11).Value = 50.761421319797
$ws.Cells.Item(16, 12).Value = 296
$ws.Cells.Item(16, 13).Value = 298.657718120805
$ws.Cells.Item(16, 14).Value = -73.901581722319
# Row 17
$ws.Cells.Item(17, 3).Value = 5
$ws.Cells.Item(17, 4).Value = 4
$ws.Cells.Item(17, 5).Value = 25
$ws.Cells.Item(17, 6).Value = 33
$ws.Cells.Item(17, 7).Value = 34
$ws.Cells.Item(17, 8).Value = -2.941176470588
$ws.Cells.Item(17, 9).Value = 430
$ws.Cells.Item(17, 10).Value = 413
$ws.Cells.Item(17, 11).Value = 4.11622276029
$ws.Cells.Item(17, 12).Value = 119.387755102041
$ws.Cells.Item(17, 13).Value = 155.952380952381
$ws.Cells.Item(17, 14).Value = -30.30794165316
# Row 18
$ws.Cells.Item(18, 3).Value = 13
$ws.Cells.Item(18, 4).Value = 10
$ws.Cells.Item(18, 5).Value = 30
$ws.Cells.Item(18, 6).Value = 50
$ws.Cells.Item(18, 7).Value = 34
$ws.Cells.Item(18, 8).Value = 47.058823529411
$ws.Cells.Item(18, 9).Value = 607
$ws.Cells.Item(18, 10).Value = 395
$ws.Cells.Item(18, 11).Value = 53.670886075949
$ws.Cells.Item(18, 12).Value = 88.509316770186
$ws.Cells.Item(18, 13).Value = 92.698412698412
$ws.Cells.Item(18, 14).Value = -75.122950819672
# Row 19
$ws.Cells.Item(19, 3).Value = 74
$ws.Cells.Item(19, 4).Value = 44
$ws.Cells.Item(19, 5).Value = 68.181818181818
$ws.Cells.Item(19, 6).Value = 211
$ws.Cells.Item(19, 7).Value = 147
$ws.Cells.Item(19, 8).Value = 43.537414965986
$ws.Cells.Item(19, 9).Value = 2085
$ws.Cells.Item(19, 10).Value = 1266
$ws.Cells.Item(19, 11).Value = 64.691943127962
$ws.Cells.Item(19, 12).Value = 96.698113207547
$ws.Cells.Item(19, 13).Value = 1.657727937591
$ws.Cells.Item(19, 14).Value = -75.653900046707
# Row 20
$ws.Cells.Item(20, 3).Value = 1
Set-TextCell 20 4 "0" 23 3
Set-TextCell 20 5 "***.*" 23 5
$ws.Cells.Item(20, 6).Value = 5
$ws.Cells.Item(20, 7).Value = 5
$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(20, 9).Value = 61
$ws.Cells.Item(20, 11).Value = 22
$ws.Cells.Item(20, 12).Value = 41.860465116279
$ws.Cells.Item(20, 13).Value = 177.272727272727
$ws.Cells.Item(20, 14).Value = -81.345565749235
# Row 21
$ws.Cells.Item(21, 3).Value = 110
$ws.Cells.Item(21, 4).Value = 70
$ws.Cells.Item(21, 5).Value = 57.142857142857
$ws.Cells.Item(21, 6).Value = 368
$ws.Cells.Item(21, 7).Value = 260
$ws.Cells.Item(21, 8).Value = 41.538461538461
$ws.Cells.Item(21, 9).Value = 3803
$ws.Cells.Item(21, 10).Value = 2538
$ws.Cells.Item(21, 11).Value = 49.842395587076
$ws.Cells.Item(21, 12).Value = 111.512791991101
$ws.Cells.Item(21, 13).Value = 39.764792355751
$ws.Cells.Item(21, 14).Value = -73.351552098661
# Row 22
$ws.Cells.Item(22, 3).Value = 3
$ws.Cells.Item(22, 4).Value = 8
$ws.Cells.Item(22, 5).Value = -62.5
$ws.Cells.Item(22, 7).Value = 16
$ws.Cells.Item(22, 8).Value = -43.75
$ws.Cells.Item(22, 9).Value = 171
$ws.Cells.Item(22, 10).Value = 142
$ws.Cells.Item(22, 11).Value = 20.422535211267
$ws.Cells.Item(22, 12).Value = 36.8
$ws.Cells.Item(22, 13).Value = 28.571428571428
# Row 24
$ws.Cells.Item(24, 3).Value = 66
$ws.Cells.Item(24, 4).Value = 62
$ws.Cells.Item(24, 5).Value = 6.451612903225
$ws.Cells.Item(24, 6).Value = 298
$ws.Cells.Item(24, 7).Value = 225
$ws.Cells.Item(24, 8).Value = 32.444444444444
$ws.Cells.Item(24, 9).Value = 3078
$ws.Cells.Item(24, 10).Value = 2039
$ws.Cells.Item(24, 11).Value = 50.956351152525
$ws.Cells.Item(24, 12).Value = 91.656288916562
$ws.Cells.Item(24, 13).Value = -29.693924166286
# Row 25
$ws.Cells.Item(25, 3).Value = 17
$ws.Cells.Item(25, 4).Value = 13
$ws.Cells.Item(25, 5).Value = 30.76923076923
$ws.Cells.Item(25, 6).Value = 75
$ws.Cells.Item(25, 7).Value = 69
$ws.Cells.Item(25, 8).Value = 8.695652173913
$ws.Cells.Item(25, 9).Value = 822
$ws.Cells.Item(25, 10).Value = 802
$ws.Cells.Item(25, 11).Value = 2.493765586034
$ws.Cells.Item(25, 12).Value = 74.52229299363
$ws.Cells.Item(25, 13).Value = 51.660516605166
# Row 26
Set-TextCell 26 3 "0" 23 3
$ws.Cells.Item(26, 6).Value = 3
Set-TextCell 26 7 "0" 23 3
Set-TextCell 26 8 "***.*" 23 5
# Row 27
$ws.Cells.Item(27, 3).Value = 4
$ws.Cells.Item(27, 4).Value = 5
$ws.Cells.Item(27, 5).Value = -20
$ws.Cells.Item(27, 6).Value = 13
$ws.Cells.Item(27, 7).Value = 15
$ws.Cells.Item(27, 8).Value = -13.333333333333
$ws.Cells.Item(27, 9).Value = 203
$ws.Cells.Item(27, 10).Value = 143
$ws.Cells.Item(27, 11).Value = 41.958041958042
$ws.Cells.Item(27, 12).Value = 86.238532110091
# Row 28
Set-TextCell 28 4 "0" 23 3
Set-TextCell 28 5 "***.*" 23 5
# Row 29
Set-TextCell 29 4 "0" 23 3
Set-TextCell 29 5 "***.*" 23 5
# Row 30
Set-TextCell 30 3 "0" 23 3
$ws.Cells.Item(30, 6).Value = 3
$ws.Cells.Item(30, 8).Value = 50
$ws.Cells.Item(30, 9).Value = 23
$ws.Cells.Item(30, 11).Value = -34.285714285714
$ws.Cells.Item(30, 12).Value = 666.666666666667
